# Natmi following Dr Hou advice -- refresh LR-pair rows 2-13 of Nlgn2-Nrxn2 sheet
# (4 sending clusters ECs/FAPs/M2/sCs x 3 target clusters ECs/FAPs/sCs)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ Row=2; A="ECs"; D="ECs"; Vals=@(3, 1, 2.580781333333333, 7.742344, 0.09250274102763278, 0.0925027410276328, 2, 0.6666666666666666, 0.09268433333333333, 0.278053, 0.1040432947149345, 0.1040432947149346, 0.2391979973591111, 2.152781976232, 0.009624289946677264, 0.009624289946677267) },
  @{ Row=3; A="ECs"; D="FAPs"; Vals=@(3, 1, 2.580781333333333, 7.742344, 0.09250274102763278, 0.0925027410276328, 3, 1, 0.5898666666666667, 1.7696, 0.6621579854471923, 0.6621579854471924, 1.522316882488889, 13.7008519424, 0.06125142864720066, 0.06125142864720069) },
  @{ Row=4; A="ECs"; D="sCs"; Vals=@(3, 1, 2.580781333333333, 7.742344, 0.09250274102763278, 0.0925027410276328, 3, 1, 0.2082736666666667, 0.624821, 0.2337987198378731, 0.2337987198378731, 0.5375087911582223, 4.837579120424, 0.02162702243375484, 0.02162702243375485) },
  @{ Row=5; A="FAPs"; D="ECs"; Vals=@(3, 1, 20.80139666666667, 62.40418999999999, 0.7455828140172008, 0.745582814017201, 2, 0.6666666666666666, 0.09268433333333333, 0.278053, 0.1040432947149345, 0.1040432947149346, 1.927963582452222, 17.35167224207, 0.07757289245318186, 0.07757289245318187) },
  @{ Row=6; A="FAPs"; D="FAPs"; Vals=@(3, 1, 20.80139666666667, 62.40418999999999, 0.7455828140172008, 0.745582814017201, 3, 1, 0.5898666666666667, 1.7696, 0.6621579854471923, 0.6621579854471924, 12.27005051377778, 110.430454624, 0.4936936141136783, 0.4936936141136785) },
  @{ Row=7; A="FAPs"; D="sCs"; Vals=@(3, 1, 20.80139666666667, 62.40418999999999, 0.7455828140172008, 0.745582814017201, 3, 1, 0.2082736666666667, 0.624821, 0.2337987198378731, 0.2337987198378731, 4.332383155554444, 38.99144839998999, 0.1743163074503405, 0.1743163074503406) },
  @{ Row=8; A="M2"; D="ECs"; Vals=@(3, 1, 0.305159, 0.915477, 0.01093778987962227, 0.01093778987962227, 2, 0.6666666666666666, 0.09268433333333333, 0.278053, 0.1040432947149345, 0.1040432947149346, 0.02828345847566667, 0.254551126281, 0.001138003695975568, 0.001138003695975568) },
  @{ Row=9; A="M2"; D="FAPs"; Vals=@(3, 1, 0.305159, 0.915477, 0.01093778987962227, 0.01093778987962227, 3, 1, 0.5898666666666667, 1.7696, 0.6621579854471923, 0.6621579854471924, 0.1800031221333333, 1.6200280992, 0.007242544911935368, 0.007242544911935371) },
  @{ Row=10; A="M2"; D="sCs"; Vals=@(3, 1, 0.305159, 0.915477, 0.01093778987962227, 0.01093778987962227, 3, 1, 0.2082736666666667, 0.624821, 0.2337987198378731, 0.2337987198378731, 0.06355658384633334, 0.5720092546169999, 0.00255724127171133, 0.00255724127171133) },
  @{ Row=11; A="sCs"; D="ECs"; Vals=@(3, 1, 4.212175, 12.636525, 0.150976655075544, 0.150976655075544, 2, 0.6666666666666666, 0.09268433333333333, 0.278053, 0.1040432947149345, 0.1040432947149346, 0.3904026317583333, 3.513623685825, 0.01570810861909984, 0.01570810861909984) },
  @{ Row=12; A="sCs"; D="FAPs"; Vals=@(3, 1, 4.212175, 12.636525, 0.150976655075544, 0.150976655075544, 3, 1, 0.5898666666666667, 1.7696, 0.6621579854471923, 0.6621579854471924, 2.484621626666667, 22.36159464, 0.09997039777437781, 0.09997039777437786) },
  @{ Row=13; A="sCs"; D="sCs"; Vals=@(3, 1, 4.212175, 12.636525, 0.150976655075544, 0.150976655075544, 3, 1, 0.2082736666666667, 0.624821, 0.2337987198378731, 0.2337987198378731, 0.8772851318916667, 7.895566187025, 0.0352981486820663, 0.03529814868206631) }
)

foreach ($rec in $rows) {
  $r = $rec.Row
  $ws.Cells.Item($r, 1).Value = $rec.A
  $ws.Cells.Item($r, 2).Value = "Nlgn2"
  $ws.Cells.Item($r, 3).Value = "Nrxn2"
  $ws.Cells.Item($r, 4).Value = $rec.D
  for ($i = 0; $i -lt $rec.Vals.Length; $i++) {
    $ws.Cells.Item($r, 5 + $i).Value = $rec.Vals[$i]
  }
}
